# Add a new employee record (MARIA LUIZA) as row 10 of the Tabela13 table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the structured table by one row (A1:P9 -> A1:P10); this also
# updates the table's <autoFilter>/ref and is how Excel grows a ListObject.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Write the new row's values in the same left-to-right order a person
# typing them into the sheet would, so new shared-string entries land in
# the same order as the target file (MARIA LUIZA, 94585301234, RUA J,
# 566, 935859334). Columns D/H/O (CPF, NUMERO, CELULAR) are switched to
# Text format before the value is entered so they are stored as text
# (matching the target, unlike the plain-number storage used by the
# pre-existing rows above).
$ws.Range("B10").Value = "MARIA LUIZA"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "94585301234"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "RUA J"

$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "566"

$ws.Range("O10").NumberFormat = "@"
$ws.Range("O10").Value = "935859334"

$ws.Range("A10").Value = 15
$ws.Range("C10").Value = "SUPORTE TÉCNICO"
$ws.Range("F10").Value = "CARLOS ALBERTO"
$ws.Range("G10").Value = "JULIA FERNANDES"
$ws.Range("J10").Value = "BAIRRO A"
$ws.Range("K10").Value = "JOÃO PESSOA"
$ws.Range("L10").Value = "PB"
$ws.Range("M10").Value = 52210901
$ws.Range("N10").Value = 83

# The whole row (like the rest of the table) is formatted as Text; apply
# this to the remaining cells too (doesn't disturb the numeric values
# already entered above). O10 additionally carries the underlined-font
# style used by several CELULAR cells in this sheet.
$ws.Range("A10:N10").NumberFormat = "@"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("O10").Font.Underline = $true

# Column O (CELULAR) did not have an explicit width before; Excel best-fit
# it to the new 9-character phone number.
$ws.Columns.Item(15).ColumnWidth = 10.1666666666667

# Leave the cursor where it ended up after entering the new row.
$ws.Range("A4").Select()
